$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.763.81'
$ws.Range('E2').Value = '  +4.02%  '
$ws.Range('D3').Value = '1.909.19'
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.23'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.695'
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '46.75'
$ws.Range('E8').Value = '  +7.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.373'
$ws.Range('E9').Value = '  +4.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.82'
$ws.Range('E10').Value = '  +6.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0758'
$ws.Range('E11').Value = '  +1.57%  '
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.73'
$ws.Range('E13').Value = '  +8.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.812'
$ws.Range('E14').Value = '  +5.01%  '
$ws.Range('D15').Value = '2.184.41'
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.09'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '1.902.78'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '36.700.86'
$ws.Range('E18').Value = '  +3.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.29'
$ws.Range('E19').Value = '  +1.25%  '
$ws.Range('D20').Value = '0.0₃0854'
$ws.Range('E20').Value = '  +3.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.67'
$ws.Range('E21').Value = '  +6.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '250.42'
$ws.Range('E22').Value = '  +2.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.14'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.54'
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.19'
$ws.Range('E26').Value = '  +1.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.23'
$ws.Range('E27').Value = '  +1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.77'
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.67'
$ws.Range('E29').Value = '  +2.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.129'
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('E31').Value = '  +7.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0611'
$ws.Range('E32').Value = '  +1.89%  '
$ws.Range('E33').Value = '  +3.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.31'
$ws.Range('E34').Value = '  +3.08%  '
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0878'
$ws.Range('E36').Value = '  +19.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.15'
$ws.Range('E37').Value = '  +59.33%  '
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.869'
$ws.Range('E39').Value = '  +1.73%  '
$ws.Range('E40').Value = '  +1.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '104.47'
$ws.Range('E41').Value = '  +7.60%  '
$ws.Range('E42').Value = '  +3.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.81'
$ws.Range('E43').Value = '  +3.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.85'
$ws.Range('E44').Value = '  +19.20%  '
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('D46').Value = '1.351.57'
$ws.Range('E46').Value = '  +3.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.38'
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0817'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.82'
$ws.Range('E49').Value = '  +3.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.41'
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('D51').Value = '2.088.78'
$ws.Range('E51').Value = '  +1.36%  '
